# Adding the changes we made on may 9th
# - Insert 9 new data rows right after the header row (rows 2-10), pushing the
#   existing 20 data rows down to rows 11-30.
# - Append 1 new data row at the very end (row 31).
# - Column A (timestamp) is simply 100*(row-2) and column B is always
#   "falling", for every data row in the sheet (1..30 -> 0,100,...,2900).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRowsTop = @(
    @(0, 'falling', -1.231432914733887, 2.921578645706177, 1.135899901390076, -0.05018257871270176, 0.03831957608461375, -0.02618168391287326),
    @(100, 'falling', -1.530599117279053, 3.089309453964233, 1.274296164512634, -0.0003359749913214444, 0.002840522676706245, 0.01902845310978596),
    @(200, 'falling', -1.078460693359375, 3.193733692169189, 1.281612634658814, 0.0197737082280218, 0.003915645778179187, 0.0108062067255378),
    @(300, 'falling', -0.9324893951416016, 2.998547315597534, 0.8108012080192566, 0.006835582219064174, 0.01575421430170528, -0.01602910399436943),
    @(400, 'falling', -1.057272911071777, 2.989111185073853, 0.8697453737258911, 0.009694431573152461, 0.04360967107117174, -0.02225992940366261),
    @(500, 'falling', -1.247255325317383, 3.032690763473511, 0.9775734543800354, -0.006389650218188861, 0.1472857224941257, -0.09550878420472185),
    @(600, 'falling', -1.120566368103027, 3.040028095245361, 0.9562293887138368, -0.0494189966470003, 0.3487124174833299, -0.1236696735024451),
    @(700, 'falling', -1.322433471679688, 3.130712985992432, 1.133776545524597, -0.1319224560260773, 0.6216225624084474, 0.00195476904511456),
    @(800, 'falling', -1.53396463394165, 3.08948278427124, 1.223363161087036, -0.2252138006687165, 0.7861163711547852, 0.1145494534075265)
)

$newRowBottom = @(2900, 'falling', -0.5697603225708008, 2.878552436828613, 1.078300476074219, 0.002003637989982902, -0.01257160693407047, -0.01087340146303169)

# Insert 9 blank rows right below the header (row 1), shifting the existing
# 20 data rows (old rows 2-21) down to rows 11-30.
$ws.Range("A2:A10").EntireRow.Insert()
# Excel's row insert copies formatting from the row above (the header) -
# the source data rows carry no explicit style, so strip it back off.
$ws.Range("A2:H10").ClearFormats()

# Fill the freshly inserted rows with the new data.
for ($i = 0; $i -lt $newRowsTop.Count; $i++) {
    $r = 2 + $i
    $rowData = $newRowsTop[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# Append the new row at the very end of the table.
$lastRow = 31
for ($c = 0; $c -lt $newRowBottom.Count; $c++) {
    $ws.Cells.Item($lastRow, $c + 1).Value = $newRowBottom[$c]
}

# The shifted (originally-existing) data rows kept their old literal
# timestamps when the rows above were inserted; re-number column A so the
# whole sheet is sequential (0, 100, 200, ... step 100) again.
for ($r = 11; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
}
